# Fruta / hortaliza, semanal
# Insert a new weekly record at row 223 of Sheet1, pushing the existing
# rows 223:235 down to 224:236.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new blank row at position 223 (shifts 223:235 -> 224:236)
$ws.Rows(223).Insert()

# Populate the new row 223 with the new weekly record
$ws.Cells.Item(223, 1).Value = 9
$ws.Cells.Item(223, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(223, 3).Value = "Metropolitana"
$ws.Cells.Item(223, 4).Value = 44753
$ws.Cells.Item(223, 5).Value = 13
$ws.Cells.Item(223, 6).Value = 100112026
$ws.Cells.Item(223, 7).Value = "Haba"
$ws.Cells.Item(223, 8).Value = "Sin especificar"
$ws.Cells.Item(223, 9).Value = "Primera"
$ws.Cells.Item(223, 10).Value = 52
$ws.Cells.Item(223, 11).Value = 22000
$ws.Cells.Item(223, 12).Value = 22000
$ws.Cells.Item(223, 13).Value = 22000
$ws.Cells.Item(223, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(223, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(223, 16).Value = 880
$ws.Cells.Item(223, 17).Value = 25
$ws.Cells.Item(223, 18).Value = "Hortaliza"
